# This edit inserts one new data row into the price list at row 69
# (pushing the existing rows 69-167 down to 70-168) and fills the new
# row with a new "Ají" (chili pepper) price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69; this shifts every row that was
# at 69..167 down by one (to 70..168) and updates the sheet dimension
# automatically (A1:R167 -> A1:R168).
$ws.Rows(69).Insert()

# Populate the newly inserted row 69 with the new record's data.
$ws.Cells.Item(69, 1).Value = 5
$ws.Cells.Item(69, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(69, 3).Value = "Maule"
$ws.Cells.Item(69, 4).Value = 44579
$ws.Cells.Item(69, 5).Value = 7
$ws.Cells.Item(69, 6).Value = 100112021
$ws.Cells.Item(69, 7).Value = "Ají"
$ws.Cells.Item(69, 8).Value = "Americana (o)"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 50
$ws.Cells.Item(69, 11).Value = 20000
$ws.Cells.Item(69, 12).Value = 20000
$ws.Cells.Item(69, 13).Value = 20000
$ws.Cells.Item(69, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(69, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(69, 16).Value = 800
$ws.Cells.Item(69, 17).Value = 25
$ws.Cells.Item(69, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date/time number format that is
# used by every other cell in column D.
$ws.Cells.Item(69, 4).NumberFormat = $ws.Cells.Item(70, 4).NumberFormat
